# Auto-generated edit script: updates market-price derived columns (H-N)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3289680.2
$ws.Range("I6").Value = 6579203
$ws.Range("J6").Value = 157.5
$ws.Range("K6").Value = 19737609
$ws.Range("L6").Value = 472.5
$ws.Range("M6").Value = -19737497
$ws.Range("N6").Value = -696.5
$ws.Range("H18").Value = 143000220
$ws.Range("I18").Value = 166833540
$ws.Range("J18").Value = 352
$ws.Range("K18").Value = 166833540
$ws.Range("L18").Value = 352
$ws.Range("M18").Value = -166833256
$ws.Range("N18").Value = -920
$ws.Range("H132").Value = 2348.6309
$ws.Range("I132").Value = 1985.1455
$ws.Range("J132").Value = 4347.8
$ws.Range("K132").Value = 5955.4365
$ws.Range("L132").Value = 13043.4
$ws.Range("M132").Value = -3425.4365
$ws.Range("N132").Value = -18103.4
$ws.Range("H135").Value = 932.3333
$ws.Range("I135").Value = 648.6667
$ws.Range("J135").Value = 2350.6667
$ws.Range("K135").Value = 5838.0003
$ws.Range("L135").Value = 21156.0003
$ws.Range("M135").Value = -3303.0003
$ws.Range("N135").Value = -26226.0003
$ws.Range("H137").Value = 3476.353
$ws.Range("I137").Value = 1375.25
$ws.Range("J137").Value = 5344
$ws.Range("K137").Value = 4125.75
$ws.Range("L137").Value = 16032
$ws.Range("M137").Value = -1575.75
$ws.Range("N137").Value = -21132
$ws.Range("H138").Value = 3259.25
$ws.Range("I138").Value = 1111.1
$ws.Range("J138").Value = 4333.325
$ws.Range("K138").Value = 3333.3
$ws.Range("L138").Value = 12999.975
$ws.Range("M138").Value = 1806.7
$ws.Range("N138").Value = -23279.975

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1694.625
$ws.Range("I61").Value = 1269.381
$ws.Range("J61").Value = 4671.3335
$ws.Range("K61").Value = 1269.381
$ws.Range("L61").Value = 4671.3335
$ws.Range("M61").Value = -1057.381
$ws.Range("N61").Value = -5095.3335
$ws.Range("H74").Value = 1464.3529
$ws.Range("I74").Value = 1379.1333
$ws.Range("J74").Value = 2103.5
$ws.Range("K74").Value = 1379.1333
$ws.Range("L74").Value = 2103.5
$ws.Range("M74").Value = -505.1333
$ws.Range("N74").Value = -3851.5
$ws.Range("H77").Value = 1464.3529
$ws.Range("I77").Value = 1379.1333
$ws.Range("J77").Value = 2103.5
$ws.Range("K77").Value = 6895.666499999999
$ws.Range("L77").Value = 10517.5
$ws.Range("M77").Value = -2527.666499999999
$ws.Range("N77").Value = -19253.5
$ws.Range("H132").Value = 2421.1785
$ws.Range("I132").Value = 2075.7646
$ws.Range("J132").Value = 2955
$ws.Range("K132").Value = 6227.293799999999
$ws.Range("L132").Value = 8865
$ws.Range("M132").Value = -3697.293799999999
$ws.Range("N132").Value = -13925
$ws.Range("H136").Value = 1694.625
$ws.Range("I136").Value = 1269.381
$ws.Range("J136").Value = 4671.3335
$ws.Range("K136").Value = 3808.143
$ws.Range("L136").Value = 14014.0005
$ws.Range("M136").Value = -1258.143
$ws.Range("N136").Value = -19114.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 13427.143
$ws.Range("I26").Value = 8798
$ws.Range("J26").Value = 25000
$ws.Range("K26").Value = 8798
$ws.Range("L26").Value = 25000
$ws.Range("M26").Value = -8506
$ws.Range("N26").Value = -25584
$ws.Range("H42").Value = 43000
$ws.Range("J42").Value = 43000
$ws.Range("L42").Value = 43000
$ws.Range("N42").Value = -43656
$ws.Range("H86").Value = 1943.7142
$ws.Range("I86").Value = 1721.2
$ws.Range("J86").Value = 2500
$ws.Range("K86").Value = 1721.2
$ws.Range("L86").Value = 2500
$ws.Range("M86").Value = -598.2
$ws.Range("N86").Value = -4746
$ws.Range("H89").Value = 1943.7142
$ws.Range("I89").Value = 1721.2
$ws.Range("J89").Value = 2500
$ws.Range("K89").Value = 8606
$ws.Range("L89").Value = 12500
$ws.Range("M89").Value = -2990
$ws.Range("N89").Value = -23732
$ws.Range("H134").Value = 2103
$ws.Range("I134").Value = 1907.45
$ws.Range("K134").Value = 5722.35
$ws.Range("M134").Value = -3187.35

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5874.75
$ws.Range("I16").Value = 3750
$ws.Range("K16").Value = 3750
$ws.Range("M16").Value = -3463
$ws.Range("H31").Value = 2314.1794
$ws.Range("I31").Value = 1715.2593
$ws.Range("J31").Value = 3661.75
$ws.Range("K31").Value = 1715.2593
$ws.Range("L31").Value = 3661.75
$ws.Range("M31").Value = -1420.2593
$ws.Range("N31").Value = -4251.75
$ws.Range("H34").Value = 2314.1794
$ws.Range("I34").Value = 1715.2593
$ws.Range("J34").Value = 3661.75
$ws.Range("K34").Value = 1715.2593
$ws.Range("L34").Value = 3661.75
$ws.Range("M34").Value = -1513.2593
$ws.Range("N34").Value = -4065.75
$ws.Range("H58").Value = 921.6863
$ws.Range("I58").Value = 749.8261
$ws.Range("K58").Value = 749.8261
$ws.Range("M58").Value = -546.8261
$ws.Range("H62").Value = 3660
$ws.Range("I62").Value = 4326.6665
$ws.Range("J62").Value = 2993.3333
$ws.Range("K62").Value = 4326.6665
$ws.Range("L62").Value = 2993.3333
$ws.Range("M62").Value = -3702.6665
$ws.Range("N62").Value = -4241.3333
$ws.Range("H65").Value = 3660
$ws.Range("I65").Value = 4326.6665
$ws.Range("J65").Value = 2993.3333
$ws.Range("K65").Value = 21633.3325
$ws.Range("L65").Value = 14966.6665
$ws.Range("M65").Value = -18513.3325
$ws.Range("N65").Value = -21206.6665
$ws.Range("H113").Value = 5874.75
$ws.Range("I113").Value = 3750
$ws.Range("K113").Value = 3750
$ws.Range("M113").Value = -1580
$ws.Range("H134").Value = 2068.138
$ws.Range("I134").Value = 1190.8462
$ws.Range("J134").Value = 9671.333000000001
$ws.Range("K134").Value = 3572.5386
$ws.Range("L134").Value = 29013.999
$ws.Range("M134").Value = -1037.5386
$ws.Range("N134").Value = -34083.999
$ws.Range("H136").Value = 921.6863
$ws.Range("I136").Value = 749.8261
$ws.Range("K136").Value = 2249.4783
$ws.Range("M136").Value = 300.5217000000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 778376.0600000001
$ws.Range("I132").Value = 1297.9
$ws.Range("K132").Value = 11681.1
$ws.Range("M132").Value = -9151.1

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 68.92856999999999
$ws.Range("I2").Value = 38.75
$ws.Range("K2").Value = 38.75
$ws.Range("M2").Value = 74.25
$ws.Range("H70").Value = 6474.577
$ws.Range("I70").Value = 7074.316
$ws.Range("J70").Value = 4846.7144
$ws.Range("K70").Value = 7074.316
$ws.Range("L70").Value = 4846.7144
$ws.Range("M70").Value = -6804.316
$ws.Range("N70").Value = -5386.7144
$ws.Range("H73").Value = 6474.577
$ws.Range("I73").Value = 7074.316
$ws.Range("J73").Value = 4846.7144
$ws.Range("K73").Value = 7074.316
$ws.Range("L73").Value = 4846.7144
$ws.Range("M73").Value = -6138.316
$ws.Range("N73").Value = -6718.7144
$ws.Range("H132").Value = 2755.3684
$ws.Range("I132").Value = 2366.3704
$ws.Range("J132").Value = 3710.182
$ws.Range("K132").Value = 7099.111199999999
$ws.Range("L132").Value = 11130.546
$ws.Range("M132").Value = -4569.111199999999
$ws.Range("N132").Value = -16190.546

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 20000
$ws.Range("J123").Value = 20000
$ws.Range("L123").Value = 20000
$ws.Range("N123").Value = -29800
$ws.Range("H132").Value = 3193.25
$ws.Range("I132").Value = 2085.1667
$ws.Range("J132").Value = 4301.3335
$ws.Range("K132").Value = 6255.500100000001
$ws.Range("L132").Value = 12904.0005
$ws.Range("M132").Value = -3725.500100000001
$ws.Range("N132").Value = -17964.0005
